$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('M3').Value = 'Printed On: 10/23/2025'
$ws.Range('B16').Value = 'Agency: KENNEDY TOWNSHIP'
$ws.Range('E21').Value = '3'
$ws.Range('I21').Value = '3'
$ws.Range('E22').Value = '3'
$ws.Range('I22').Value = '3'
$ws.Range('E24').Value = '3'
$ws.Range('I24').Value = '3'
$ws.Range('E25').Value = '2'
$ws.Range('I25').Value = '2'
$ws.Range('E29').Value = '4'
$ws.Range('I29').Value = '4'
$ws.Range('E30').Value = '0'
$ws.Range('I30').Value = '0'
$ws.Range('E33').Value = '4'
$ws.Range('I33').Value = '4'
$ws.Range('E34').Value = '2'
$ws.Range('I34').Value = '2'
$ws.Range('J34').Value = '1'
$ws.Range('E35').Value = '1'
$ws.Range('I35').Value = '1'
$ws.Range('J35').Value = '1'
$ws.Range('E36').Value = '1'
$ws.Range('I36').Value = '1'
$ws.Range('E38').Value = '92'
$ws.Range('I38').Value = '92'
$ws.Range('J38').Value = '32'
$ws.Range('N38').Value = '2'
$ws.Range('E39').Value = '7'
$ws.Range('I39').Value = '7'
$ws.Range('J39').Value = '1'
$ws.Range('E40').Value = '6'
$ws.Range('I40').Value = '6'
$ws.Range('J40').Value = '1'
$ws.Range('E42').Value = '1'
$ws.Range('I42').Value = '1'
$ws.Range('E45').Value = '111'
$ws.Range('I45').Value = '111'
$ws.Range('J45').Value = '37'
$ws.Range('N45').Value = '2'
$ws.Range('E46').Value = '31'
$ws.Range('I46').Value = '31'
$ws.Range('J46').Value = '13'
$ws.Range('N46').Value = '4'
$ws.Range('E47').Value = '1'
$ws.Range('I47').Value = '1'
$ws.Range('E48').Value = '38'
$ws.Range('I48').Value = '38'
$ws.Range('J48').Value = '5'
$ws.Range('E50').Value = '1'
$ws.Range('I50').Value = '1'
$ws.Range('J50').Value = '1'
$ws.Range('N50').Value = '0'
$ws.Range('E51').Value = '17'
$ws.Range('I51').Value = '17'
$ws.Range('N51').Value = '0'
$ws.Range('E52').Value = '1'
$ws.Range('I52').Value = '1'
$ws.Range('J52').Value = '1'
$ws.Range('J54').Value = '0'
$ws.Range('E55').Value = '6'
$ws.Range('I55').Value = '6'
$ws.Range('J55').Value = '5'
$ws.Range('E56').Value = '0'
$ws.Range('I56').Value = '0'
$ws.Range('J56').Value = '0'
$ws.Range('E57').Value = '0'
$ws.Range('I57').Value = '0'
$ws.Range('J57').Value = '0'
$ws.Range('E61').Value = '6'
$ws.Range('I61').Value = '6'
$ws.Range('J61').Value = '5'
$ws.Range('E62').Value = '4'
$ws.Range('I62').Value = '4'
$ws.Range('J62').Value = '4'
$ws.Range('E63').Value = '1'
$ws.Range('I63').Value = '1'
$ws.Range('J63').Value = '0'
$ws.Range('E64').Value = '0'
$ws.Range('I64').Value = '0'
$ws.Range('E65').Value = '1'
$ws.Range('I65').Value = '1'
$ws.Range('E70').Value = '0'
$ws.Range('I70').Value = '0'
$ws.Range('E71').Value = '4'
$ws.Range('I71').Value = '4'
$ws.Range('J71').Value = '2'
$ws.Range('E72').Value = '1'
$ws.Range('I72').Value = '1'
$ws.Range('E73').Value = '5'
$ws.Range('I73').Value = '5'
$ws.Range('J73').Value = '4'
$ws.Range('E74').Value = '19'
$ws.Range('I74').Value = '19'
$ws.Range('J74').Value = '5'
$ws.Range('E76').Value = '43'
$ws.Range('I76').Value = '43'
$ws.Range('J76').Value = '10'
$ws.Range('N76').Value = '2'
$ws.Range('E77').Value = '169'
$ws.Range('I77').Value = '169'
$ws.Range('J77').Value = '48'
$ws.Range('N77').Value = '6'
